# Refresh cryptos list (GitHub Actions scheduled update).
# Price (col D) and Volume(1h) (col E) are plain text cells in the source
# sheet (t="inlineStr"), not numbers, so numeric-looking prices are forced
# back to Text via NumberFormat "@" before assignment to avoid Excel
# auto-converting them to the Number type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.599.52"
$ws.Range("E2").Value = "  -0.49%  "

$ws.Range("D3").Value = "2.289.21"
$ws.Range("E3").Value = "  -0.76%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "96.05"
$ws.Range("E5").Value = "  +2.57%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "267.93"
$ws.Range("E6").Value = "  -1.19%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.616"
$ws.Range("E7").Value = "  -2.19%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("E9").Value = "  -2.14%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.60"
$ws.Range("E10").Value = "  +1.89%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0935"
$ws.Range("E11").Value = "  -0.23%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.99"
$ws.Range("E12").Value = "  -1.31%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.105"
$ws.Range("E13").Value = "  +0.21%  "

$ws.Range("D14").Value = "2.632.17"
$ws.Range("E14").Value = "  -0.82%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.39"
$ws.Range("E15").Value = "  +0.38%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.844"
$ws.Range("E16").Value = "  -0.80%  "

$ws.Range("D17").Value = "2.289.61"
$ws.Range("E17").Value = "  -0.78%  "

$ws.Range("D18").Value = "43.572.63"
$ws.Range("E18").Value = "  -0.42%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000109"
$ws.Range("E19").Value = "  +1.91%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.19"
$ws.Range("E20").Value = "  -1.04%  "

$ws.Range("E21").Value = "  +1.20%  "

$ws.Range("E22").Value = "  +11.78%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.95"
$ws.Range("E23").Value = "  -3.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.11"
$ws.Range("E24").Value = "  -6.67%  "

$ws.Range("E25").Value = "  -0.08%  "

$ws.Range("E26").Value = "  +1.80%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.20"
$ws.Range("E27").Value = "  -1.78%  "

$ws.Range("E28").Value = "  +2.97%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.11"
$ws.Range("E29").Value = "  +2.53%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.23"
$ws.Range("E30").Value = "  -6.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.54"
$ws.Range("E31").Value = "  +2.26%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.14"
$ws.Range("E32").Value = "  -1.48%  "

$ws.Range("E33").Value = "  -0.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.35"
$ws.Range("E34").Value = "  -4.28%  "

$ws.Range("E35").Value = "  -0.44%  "

$ws.Range("E36").Value = "  -2.39%  "

$ws.Range("E37").Value = "  -0.22%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.37"
$ws.Range("E38").Value = "  -3.00%  "

$ws.Range("E39").Value = "  -2.76%  "

$ws.Range("E40").Value = "  +2.67%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.31"
$ws.Range("E41").Value = "  +0.96%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.30"
$ws.Range("E42").Value = "  +0.82%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.65"
$ws.Range("E43").Value = "  +6.40%  "

$ws.Range("E44").Value = "  +2.64%  "

$ws.Range("E45").Value = "  -1.24%  "

$ws.Range("E46").Value = "  -4.86%  "

$ws.Range("E47").Value = "  -1.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "97.13"
$ws.Range("E48").Value = "  -3.21%  "

$ws.Range("E49").Value = "  -0.91%  "

$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.51"
$ws.Range("E50").Value = "  +9.36%  "

$ws.Range("B51").Value = "TheGraph"
$ws.Range("C51").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.185"
$ws.Range("E51").Value = "  +6.61%  "
